$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5; existing rows 5-17 shift down to 6-18
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new weekly price record
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C5").Value = "Los Lagos"
$ws.Range("D5").Value = 44792
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 100112035
$ws.Range("G5").Value = "Bruselas (repollito)"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 24000
$ws.Range("L5").Value = 24000
$ws.Range("M5").Value = 24000
$ws.Range("N5").Value = "$/malla 15 kilos"
$ws.Range("O5").Value = "Provincia de Quillota"
$ws.Range("P5").Value = 1600
$ws.Range("Q5").Value = 15
$ws.Range("R5").Value = "Hortaliza"
